# Remove the <w:contextualSpacing w:val="0"/> paragraph-property element
# from every paragraph in the document (it is present, identically, on
# every paragraph's pPr). There is no COM/VBA-exposed
# ParagraphFormat.ContextualSpacing property in this object model, so we
# round-trip each paragraph's raw WordprocessingML through
# Range.WordOpenXML / Range.InsertXML, stripping the element out of the
# markup before writing it back onto the exact same range.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $rng = $para.Range
    $xml = $rng.WordOpenXML

    if ($xml -like "*<w:contextualSpacing*") {
        $newXml = [System.Text.RegularExpressions.Regex]::Replace(
            $xml, '<w:contextualSpacing\b[^/]*/>', ''
        )
        [void]$rng.InsertXML($newXml)
    }
}

Write-Output "contextualSpacing removed from $count paragraphs"
